$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (old D "Logs" / E "Bonus Xp" shift right to E / F)
$ws.Columns.Item(4).EntireColumn.Insert()

# The new column D comes in formatted like the old column C/D (style s="1");
# the new "Cook Time" data column uses the default (unstyled) cells, so strip
# the inherited formatting from the data rows before writing values.
$ws.Range("D2:D18").ClearFormats()

# Header + new "Cook Time" column data
$ws.Range("D1").Value2 = "Cook Time"
$ws.Range("D2").Value2 = 2
$ws.Range("D3").Value2 = 3
$ws.Range("D4").Value2 = 2
$ws.Range("D5").Value2 = 3
$ws.Range("D6").Value2 = 4
$ws.Range("D7").Value2 = 4
$ws.Range("D8").Value2 = 4
$ws.Range("D9").Value2 = 4
$ws.Range("D10").Value2 = 5
$ws.Range("D11").Value2 = 6
$ws.Range("D12").Value2 = 6
$ws.Range("D13").Value2 = 7
$ws.Range("D14").Value2 = 7
$ws.Range("D15").Value2 = 8
$ws.Range("D16").Value2 = 9
$ws.Range("D17").Value2 = 10
$ws.Range("D18").Value2 = 11

# Empty styled cell at G1 (matches the surrounding header formatting)
$ws.Range("G1").Font.Size = 11

# Update the sheet selection to match the authored state
$ws.Range("G1:G18").Select()
